$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new columns before the old "Terms Typically Offered" column (D),
# which shifts it to column G and makes room for the new Corequisites (D),
# Concurrent (E) and Recommended (F) columns.
$ws.Columns("D:F").Insert()

# Header row
$ws.Range("D1").Value = "Corequisites"
$ws.Range("E1").Value = "Concurrent"
$ws.Range("F1").Value = "Recommended"

# Default every data row's new columns to "NA"; specific rows are
# overridden below with values split out of the old Prerequisites text.
$ws.Range("D2:F104").Value = "NA"

$exceptions = @{}
$exceptions[11] = @{ C = "CE 204 or CE 208 (CE 208 may be taken concurrently)."; D = "CE 113."; E = "NA"; F = "NA"; G = "F, W, SP " }
$exceptions[14] = @{ C = "CE 222."; D = "CE 321."; E = "NA"; F = "NA"; G = "F, W, SP " }
$exceptions[15] = @{ C = "ME 341 or ENVE 264."; D = "NA"; E = "CE 337."; F = "NA"; G = "F, W, SP " }
$exceptions[16] = @{ C = "ME 341 or ENVE 264."; D = "NA"; E = "CE 336."; F = "NA"; G = "F, W, SP " }
$exceptions[17] = @{ C = "CE 207 or CE 208."; D = "CE 251."; E = "NA"; F = "NA"; G = "F, W, SP " }
$exceptions[21] = @{ C = "CE 207 or CE 208; ME 341 or ENVE 264."; D = "NA"; E = "CE 382 (CE majors only)."; F = "NA"; G = "F, W, SP " }
$exceptions[22] = @{ C = "NA"; D = "CE 381."; E = "NA"; F = "NA"; G = "F, W, SP" }
$exceptions[44] = @{ C = "CE 355."; D = "CE 356."; E = "NA"; F = "NA"; G = "SP " }
$exceptions[45] = @{ C = "CE 355."; D = "NA"; E = "CE 356."; F = "NA"; G = "TBD " }
$exceptions[78] = @{ C = "CE 434."; D = "ENVE 331."; E = "NA"; F = "NA"; G = "F " }
$exceptions[81] = @{ C = "CE 454."; D = "NA"; E = "NA"; F = "Concurrent enrollment in CE 557."; G = "SP " }
$exceptions[82] = @{ C = "CE 356 and senior or graduate standing."; D = "NA"; E = "NA"; F = "CE 454 and CE 407."; G = "W " }
$exceptions[103] = @{ C = "Graduate standing."; D = "NA"; E = "NA"; F = "Student should be in the final quarter of completing graduate coursework (45 units of 400 and 500 level coursework) and prepared to take the MS exam."; G = "F, W, SP " }

foreach ($rowNum in $exceptions.Keys) {
    $vals = $exceptions[$rowNum]
    $ws.Cells.Item($rowNum, 3).Value = $vals.C
    $ws.Cells.Item($rowNum, 4).Value = $vals.D
    $ws.Cells.Item($rowNum, 5).Value = $vals.E
    $ws.Cells.Item($rowNum, 6).Value = $vals.F
    $ws.Cells.Item($rowNum, 7).Value = $vals.G
}

Write-Output "edit complete"
